$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Bvals = @(0.8447786353271454,0.7469369850426801,0.6868693826353081,0.6623939359138831,0.6583299836951539,0.6865392867911453,0.8110416393366791,1.055233440527388,1.234664091416676,1.316297573398799,1.347211107082046,1.340553303561023,1.318840838045048,1.305541398331229,1.229329286213442,1.182577592394523,1.155688195025334,1.146584089435976,1.187554297626548,1.325218301574807,1.415193503754551,1.367171911917694,1.185304364116803,0.98916890488033)
for ($i = 0; $i -lt $Bvals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $Bvals[$i]
}

$Cvals = @(0.1861055023120173,0.1632408360593445,0.1491385002046854,0.1433759593091395,0.1424181501443513,0.1490608479560649,0.1782350242866926,0.2349369906712298,0.2762822943695085,0.2950228324378088,0.3021095396570104,0.3005837351251728,0.295606060199276,0.2925557918462403,0.2750561776374525,0.2643032605548683,0.2581121112973221,0.2560148091232861,0.2654485862025524,0.2970683949890542,0.3176758123416903,0.306682613366462,0.2649308129821577,0.2196524437648861)
for ($i = 0; $i -lt $Cvals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $Cvals[$i]
}

$Dvals = @(0.0790056553439058,0.0715886933525951,0.0670720335749877,0.06524078403315059,0.06493726754327156,0.06704729903968598,0.07644048017949956,0.09516147321394897,0.1091067444823608,0.1154939721189976,0.1179190041329292,0.1173964474182156,0.1156933536947662,0.1146509859978835,0.1086902052685446,0.1050446448827103,0.1029519002954657,0.1022440321986267,0.1054322976249438,0.1161934212457822,0.123263356372135,0.1194865946258261,0.1052570300589935,0.09006395580350102)
for ($i = 0; $i -lt $Dvals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 4).Value = $Dvals[$i]
}

$Evals = @(0.1070133005347813,0.1015708116470222,0.09833168479296006,0.0970371897197353,0.09682376911598567,0.09831412405240059,0.1051152445365986,0.1192813758646452,0.1302168616632216,0.1353111270698832,0.1372577730138147,0.1368377414523039,0.1354709248889847,0.1346360073459749,0.1298863772595737,0.1270035297142798,0.125356633689016,0.1248009430797339,0.1273092477529261,0.1358719127083035,0.1415706166682824,0.1385196149155803,0.1271709999777997,0.1153580129280201)
for ($i = 0; $i -lt $Evals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $Evals[$i]
}

$Gvals = @(0.002412298101118912,0.002415983506122934,0.002418365445566493,0.002419366147272076,0.002419534130498324,0.002418378819589132,0.002413544170951445,0.002405003900566905,0.002399296475996046,0.002396821837240383,0.00239590215523966,0.002396099452174024,0.002396745826038944,0.002397144013226277,0.002399460640483165,0.002400912921821061,0.002401759695658964,0.002402048369615046,0.002400757138630222,0.00239655549845061,0.002393910920298274,0.002395313129699665,0.002400827531264919,0.002407214232089262)
for ($i = 0; $i -lt $Gvals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 7).Value = $Gvals[$i]
}

$Mvals = @(0.3492223880770027,0.3124482855449244,0.2899595116798181,0.2808177462020822,0.2793011216366779,0.2898361314247779,0.336523719514048,0.4288118281979081,0.4970905578935856,0.5282618005665398,0.5400818629507285,0.5375354770814909,0.5292339183454118,0.5241510895333334,0.4950556960225896,0.477235189266338,0.4669957292259426,0.4635306071228342,0.479131132059905,0.53167184604537,0.5661047210531081,0.5477185208514044,0.4782739583264259,0.4037640928077764)
for ($i = 0; $i -lt $Mvals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 13).Value = $Mvals[$i]
}

$Ovals = @(2.559241357825584,2.540367486991528,2.531097262241502,2.527899745928153,2.527403759838535,2.531051794047016,2.552250758640952,2.612350048252438,2.667992631489426,2.695842307814928,2.706756392010107,2.70438943140573,2.696732822633578,2.692090947804445,2.666223918429722,2.651007484195787,2.642494056645546,2.63965245219498,2.652602576320874,2.698971741517198,2.731423202168912,2.713905808238053,2.651880704292978,2.59408816606134)
for ($i = 0; $i -lt $Ovals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 15).Value = $Ovals[$i]
}
